$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume data.
# D-column price cells are forced to Text format ("@") before assignment
# so Excel does not auto-convert numeric-looking strings (e.g. "239.29")
# into real numbers, which would silently drop formatting like trailing
# zeros ("74.30" -> 74.3) and change the stored cell type.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.702.45"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.348.81"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.29"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.668"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.30"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.86"
$ws.Range("E11").Value = "  +4.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "37.14"
$ws.Range("E12").Value = "  +15.79%  "
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.701.72"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.34"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.925"
$ws.Range("E17").Value = "  +4.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.355.05"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.654.16"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000103"
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.99"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.97"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  +3.40%  "
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.66"
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.75"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.133"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0754"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.52"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("E38").Value = "  +6.58%  "
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0279"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.50"
$ws.Range("E41").Value = "  +17.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.08"
$ws.Range("E42").Value = "  +13.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.82"
$ws.Range("E43").Value = "  +10.48%  "
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("E47").Value = "  +3.79%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "97.97"
$ws.Range("E51").Value = "  -0.95%  "
